# Update cryptocurrency price/volume figures to the latest scraped snapshot.
# Each target cell is stored as plain text (not a number/percentage), so we
# force text entry via a leading apostrophe (Formula = "'<text>") and then
# strip the resulting Quote-Prefix formatting with ClearFormats() so the cell's
# style index is left exactly as it was before the edit (style 0 / no style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $cellRef, $text) {
    $cell = $sheet.Range($cellRef)
    $cell.Formula = "'" + $text
    $cell.ClearFormats()
}

Set-TextValue $ws "D2" "263.15"
Set-TextValue $ws "E2" "0.95%"
Set-TextValue $ws "D3" "26.65"
Set-TextValue $ws "E3" "-1.86%"
Set-TextValue $ws "D4" "4.679"
Set-TextValue $ws "E4" "0.11%"
Set-TextValue $ws "D5" "0.06107"
Set-TextValue $ws "E5" "-1.12%"
Set-TextValue $ws "E6" "0.50%"
Set-TextValue $ws "D7" "0.8499"
Set-TextValue $ws "E7" "-0.17%"
Set-TextValue $ws "D8" "0.9114"
Set-TextValue $ws "E8" "-0.61%"
Set-TextValue $ws "D9" "0.04811"
Set-TextValue $ws "E9" "0.26%"
Set-TextValue $ws "D10" "0.07097"
Set-TextValue $ws "E10" "0.22%"
Set-TextValue $ws "D11" "0.03129"
Set-TextValue $ws "E11" "0.39%"
Set-TextValue $ws "D12" "0.09043"
Set-TextValue $ws "D13" "0.001535"
Set-TextValue $ws "E13" "-0.30%"
Set-TextValue $ws "D14" "0.0006201"
Set-TextValue $ws "E14" "0.61%"
Set-TextValue $ws "D15" "0.005977"
Set-TextValue $ws "E15" "-2.99%"
Set-TextValue $ws "D16" "3.456"
Set-TextValue $ws "E16" "0.19%"
Set-TextValue $ws "D17" "3.165"
Set-TextValue $ws "E17" "0.34%"
Set-TextValue $ws "E19" "-0.19%"
Set-TextValue $ws "D20" "0.1411"
Set-TextValue $ws "E20" "-0.01%"
Set-TextValue $ws "E21" "-1.45%"
Set-TextValue $ws "D22" "4.111"
Set-TextValue $ws "E22" "0.68%"
Set-TextValue $ws "D23" "0.04240"
Set-TextValue $ws "E23" "-0.04%"
Set-TextValue $ws "D24" "0.001180"
Set-TextValue $ws "E24" "-2.95%"
Set-TextValue $ws "D25" "0.004064"
Set-TextValue $ws "E25" "6.89%"
Set-TextValue $ws "E26" "0.00%"
Set-TextValue $ws "E27" "23.05%"
Set-TextValue $ws "D40" "0.03927"
Set-TextValue $ws "E40" "1.29%"
Set-TextValue $ws "E41" "0.05%"
Set-TextValue $ws "D42" "0.004172"
Set-TextValue $ws "E42" "1.99%"
Set-TextValue $ws "E43" "-3.94%"
Set-TextValue $ws "E44" "-28.89%"
Set-TextValue $ws "D45" "0.00005091"
Set-TextValue $ws "E45" "-1.19%"
Set-TextValue $ws "E46" "0.00%"
Set-TextValue $ws "D48" "0.2583"
Set-TextValue $ws "E48" "59.68%"
Set-TextValue $ws "E49" "0.00%"
Set-TextValue $ws "E50" "0.00%"
